$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "My Series" to "Data"
$ws.Name = "Data"

# Reorder the three data rows (2, 3, 4): the row that used to be 4
# ("Value Added in Industry") moves up to row 2, the old row 2
# ("Govt Revenue - Tax ; Individual Income") moves down to row 3, and the
# old row 3 ("Govt Revenue") moves down to row 4. Capture all three rows
# first so the shuffle doesn't clobber data it still needs to read.
$row2 = $ws.Range("A2:U2").Value()
$row3 = $ws.Range("A3:U3").Value()
$row4 = $ws.Range("A4:U4").Value()

$ws.Range("A2:U2").Value = $row4
$ws.Range("A3:U3").Value = $row2
$ws.Range("A4:U4").Value = $row3
